$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 22.06580366666666
$ws.Range("H2").Value = 66.19741099999999
$ws.Range("I2").Value = 0.1403713626377477
$ws.Range("J2").Value = 0.1403713626377477
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.360972666666666
$ws.Range("N2").Value = 19.082918
$ws.Range("O2").Value = 0.03115862208643261
$ws.Range("P2").Value = 0.03115862208643262
$ws.Range("Q2").Value = 140.3599739916997
$ws.Range("R2").Value = 1263.239765925298
$ws.Range("S2").Value = 0.004373778240187169
$ws.Range("T2").Value = 0.004373778240187169

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 22.06580366666666
$ws.Range("H3").Value = 66.19741099999999
$ws.Range("I3").Value = 0.1403713626377477
$ws.Range("J3").Value = 0.1403713626377477
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 107.3681206666667
$ws.Range("N3").Value = 322.104362
$ws.Range("O3").Value = 0.5259325690101214
$ws.Range("P3").Value = 0.5259325690101214
$ws.Range("Q3").Value = 2369.163870689642
$ws.Range("R3").Value = 21322.47483620678
$ws.Range("S3").Value = 0.07382587136752204
$ws.Range("T3").Value = 0.07382587136752203

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 22.06580366666666
$ws.Range("H4").Value = 66.19741099999999
$ws.Range("I4").Value = 0.1403713626377477
$ws.Range("J4").Value = 0.1403713626377477
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 59.10257333333334
$ws.Range("N4").Value = 177.30772
$ws.Range("O4").Value = 0.2895083571855736
$ws.Range("P4").Value = 0.2895083571855736
$ws.Range("Q4").Value = 1304.145779368102
$ws.Range("R4").Value = 11737.31201431292
$ws.Range("S4").Value = 0.04063868259315475
$ws.Range("T4").Value = 0.04063868259315474

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 22.06580366666666
$ws.Range("H5").Value = 66.19741099999999
$ws.Range("I5").Value = 0.1403713626377477
$ws.Range("J5").Value = 0.1403713626377477
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 31.31640666666667
$ws.Range("N5").Value = 93.94922000000001
$ws.Range("O5").Value = 0.1534004517178723
$ws.Range("P5").Value = 0.1534004517178724
$ws.Range("Q5").Value = 691.0216810521578
$ws.Range("R5").Value = 6219.19512946942
$ws.Range("S5").Value = 0.02153303043688377
$ws.Range("T5").Value = 0.02153303043688377

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 98.09611
$ws.Range("H6").Value = 294.28833
$ws.Range("I6").Value = 0.6240373039738243
$ws.Range("J6").Value = 0.6240373039738243
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 6.360972666666666
$ws.Range("N6").Value = 19.082918
$ws.Range("O6").Value = 0.03115862208643261
$ws.Range("P6").Value = 0.03115862208643262
$ws.Range("Q6").Value = 623.9866744163265
$ws.Range("R6").Value = 5615.88006974694
$ws.Range("S6").Value = 0.01944414252235667
$ws.Range("T6").Value = 0.01944414252235667

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 98.09611
$ws.Range("H7").Value = 294.28833
$ws.Range("I7").Value = 0.6240373039738243
$ws.Range("J7").Value = 0.6240373039738243
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 107.3681206666667
$ws.Range("N7").Value = 322.104362
$ws.Range("O7").Value = 0.5259325690101214
$ws.Range("P7").Value = 0.5259325690101214
$ws.Range("Q7").Value = 10532.39497541061
$ws.Range("R7").Value = 94791.55477869546
$ws.Range("S7").Value = 0.3282015424371035
$ws.Range("T7").Value = 0.3282015424371035

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 98.09611
$ws.Range("H8").Value = 294.28833
$ws.Range("I8").Value = 0.6240373039738243
$ws.Range("J8").Value = 0.6240373039738243
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 59.10257333333334
$ws.Range("N8").Value = 177.30772
$ws.Range("O8").Value = 0.2895083571855736
$ws.Range("P8").Value = 0.2895083571855736
$ws.Range("Q8").Value = 5797.732534989734
$ws.Range("R8").Value = 52179.5928149076
$ws.Range("S8").Value = 0.1806640146959763
$ws.Range("T8").Value = 0.1806640146959763

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 98.09611
$ws.Range("H9").Value = 294.28833
$ws.Range("I9").Value = 0.6240373039738243
$ws.Range("J9").Value = 0.6240373039738243
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 31.31640666666667
$ws.Range("N9").Value = 93.94922000000001
$ws.Range("O9").Value = 0.1534004517178723
$ws.Range("P9").Value = 0.1534004517178724
$ws.Range("Q9").Value = 3072.017673178067
$ws.Range("R9").Value = 27648.1590586026
$ws.Range("S9").Value = 0.09572760431838785
$ws.Range("T9").Value = 0.09572760431838788

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 29.56610533333334
$ws.Range("H10").Value = 88.69831600000001
$ws.Range("I10").Value = 0.1880844476016372
$ws.Range("J10").Value = 0.1880844476016372
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 6.360972666666666
$ws.Range("N10").Value = 19.082918
$ws.Range("O10").Value = 0.03115862208643261
$ws.Range("P10").Value = 0.03115862208643262
$ws.Range("Q10").Value = 188.0691878851209
$ws.Range("R10").Value = 1692.622690966088
$ws.Range("S10").Value = 0.005860452223154852
$ws.Range("T10").Value = 0.005860452223154853

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 29.56610533333334
$ws.Range("H11").Value = 88.69831600000001
$ws.Range("I11").Value = 0.1880844476016372
$ws.Range("J11").Value = 0.1880844476016372
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 107.3681206666667
$ws.Range("N11").Value = 322.104362
$ws.Range("O11").Value = 0.5259325690101214
$ws.Range("P11").Value = 0.5259325690101214
$ws.Range("Q11").Value = 3174.457165072711
$ws.Range("R11").Value = 28570.1144856544
$ws.Range("S11").Value = 0.09891973671797864
$ws.Range("T11").Value = 0.09891973671797864

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 29.56610533333334
$ws.Range("H12").Value = 88.69831600000001
$ws.Range("I12").Value = 0.1880844476016372
$ws.Range("J12").Value = 0.1880844476016372
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 59.10257333333334
$ws.Range("N12").Value = 177.30772
$ws.Range("O12").Value = 0.2895083571855736
$ws.Range("P12").Value = 0.2895083571855736
$ws.Range("Q12").Value = 1747.432908644392
$ws.Range("R12").Value = 15726.89617779952
$ws.Range("S12").Value = 0.05445201943730608
$ws.Range("T12").Value = 0.05445201943730608

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 29.56610533333334
$ws.Range("H13").Value = 88.69831600000001
$ws.Range("I13").Value = 0.1880844476016372
$ws.Range("J13").Value = 0.1880844476016372
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 31.31640666666667
$ws.Range("N13").Value = 93.94922000000001
$ws.Range("O13").Value = 0.1534004517178723
$ws.Range("P13").Value = 0.1534004517178724
$ws.Range("Q13").Value = 925.9041781681691
$ws.Range("R13").Value = 8333.137603513522
$ws.Range("S13").Value = 0.02885223922319764
$ws.Range("T13").Value = 0.02885223922319765

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 7.467887999999999
$ws.Range("H14").Value = 22.403664
$ws.Range("I14").Value = 0.04750688578679087
$ws.Range("J14").Value = 0.04750688578679088
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 6.360972666666666
$ws.Range("N14").Value = 19.082918
$ws.Range("O14").Value = 0.03115862208643261
$ws.Range("P14").Value = 0.03115862208643262
$ws.Range("Q14").Value = 47.50303144572799
$ws.Range("R14").Value = 427.527283011552
$ws.Range("S14").Value = 0.001480249100733934
$ws.Range("T14").Value = 0.001480249100733934

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 7.467887999999999
$ws.Range("H15").Value = 22.403664
$ws.Range("I15").Value = 0.04750688578679087
$ws.Range("J15").Value = 0.04750688578679088
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 107.3681206666667
$ws.Range("N15").Value = 322.104362
$ws.Range("O15").Value = 0.5259325690101214
$ws.Range("P15").Value = 0.5259325690101214
$ws.Range("Q15").Value = 801.813099909152
$ws.Range("R15").Value = 7216.317899182369
$ws.Range("S15").Value = 0.02498541848751735
$ws.Range("T15").Value = 0.02498541848751735

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 7.467887999999999
$ws.Range("H16").Value = 22.403664
$ws.Range("I16").Value = 0.04750688578679087
$ws.Range("J16").Value = 0.04750688578679088
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 59.10257333333334
$ws.Range("N16").Value = 177.30772
$ws.Range("O16").Value = 0.2895083571855736
$ws.Range("P16").Value = 0.2895083571855736
$ws.Range("Q16").Value = 441.37139816512
$ws.Range("R16").Value = 3972.34258348608
$ws.Range("S16").Value = 0.0137536404591365
$ws.Range("T16").Value = 0.0137536404591365

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 7.467887999999999
$ws.Range("H17").Value = 22.403664
$ws.Range("I17").Value = 0.04750688578679087
$ws.Range("J17").Value = 0.04750688578679088
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 31.31640666666667
$ws.Range("N17").Value = 93.94922000000001
$ws.Range("O17").Value = 0.1534004517178723
$ws.Range("P17").Value = 0.1534004517178724
$ws.Range("Q17").Value = 233.86741754912
$ws.Range("R17").Value = 2104.80675794208
$ws.Range("S17").Value = 0.007287577739403089
$ws.Range("T17").Value = 0.007287577739403092

Write-Output "Applied Natmi Gas6-Axl updates per Dr Hou advice"
